# Update the title-slide presenter/date line
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(5).TextFrame.TextRange.Text = "Presenter Name | November 15, 2025"

# Bold the last row of each summary table across the affected slides
$tableSlides = @(4, 8, 12, 17, 21)

foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(3)
    $tbl = $shape.Table
    $lastRow = $tbl.Rows.Count
    $colCount = $tbl.Columns.Count
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $tbl.Cell($lastRow, $c)
        $cell.Shape.TextFrame.TextRange.Font.Bold = 1
    }
}
